# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped values, as produced by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Some of the new Price values look like plain decimal numbers (e.g. "534.75").
# The source column stores prices as text, so force those specific cells to a
# text format while assigning them (then restore the default "Normal" style)
# to stop Excel from silently re-typing them as numeric values.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "59.055.61"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.503.44"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.28%  "
Set-TextValue "D5" "534.75"
$ws.Range("E5").Value = "  +2.75%  "
Set-TextValue "D6" "134.22"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "2.505.39"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "2.946.04"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "58.811.77"
$ws.Range("E15").Value = "  +0.81%  "
Set-TextValue "D16" "22.33"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "2.506.69"
$ws.Range("E18").Value = "  -0.48%  "
Set-TextValue "D19" "10.63"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  +1.43%  "
Set-TextValue "D21" "321.13"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  -0.20%  "
Set-TextValue "D24" "65.95"
$ws.Range("E24").Value = "  +3.45%  "
Set-TextValue "D25" "0.410"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  +0.64%  "
Set-TextValue "D30" "171.92"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  +0.08%  "
Set-TextValue "D35" "0.998"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +3.26%  "
Set-TextValue "D40" "0.829"
$ws.Range("E40").Value = "  +5.91%  "
Set-TextValue "D41" "36.46"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  +0.90%  "
Set-TextValue "D43" "275.08"
$ws.Range("E43").Value = "  -1.89%  "
Set-TextValue "D44" "131.12"
$ws.Range("E44").Value = "  +6.67%  "
Set-TextValue "D45" "5.02"
$ws.Range("E45").Value = "  -1.39%  "
Set-TextValue "D46" "0.591"
$ws.Range("E46").Value = "  -1.42%  "
Set-TextValue "D47" "0.0936"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "1.747.52"
$ws.Range("E51").Value = "  -0.11%  "
